# Fruta / hortaliza, semanal
# Insert 4 new daily price records (rows 302-305) for "Murcott" mandarinas,
# pushing the existing rows 302-345 down to 306-349.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows above the current row 302; this shifts rows 302:345
# down to 306:349 and grows the sheet dimension to A1:T349.
$ws.Rows("302:305").Insert()

# Columns that are constant across every data row in this sheet.
$mercadoId   = 10
$mercado     = "Vega Modelo de Temuco"
$region      = "La Araucanía"
$codreg      = 9
$tipo        = "Fruta"
$productoId  = 100102
$producto    = "Cítricos"
$categoriaId = 100102004
$categoria   = "Mandarina"
$origen      = "Región de O'Higgins"

# --- Row 302 ---
$ws.Range("A302").Value = $mercadoId
$ws.Range("B302").Value = $mercado
$ws.Range("C302").Value = $region
$ws.Range("D302").Value = 44449
$ws.Range("E302").Value = $codreg
$ws.Range("F302").Value = $tipo
$ws.Range("G302").Value = $productoId
$ws.Range("H302").Value = $producto
$ws.Range("I302").Value = $categoriaId
$ws.Range("J302").Value = $categoria
$ws.Range("K302").Value = "Murcott"
$ws.Range("L302").Value = "Primera"
$ws.Range("M302").Value = 200
$ws.Range("N302").Value = 7000
$ws.Range("O302").Value = 7000
$ws.Range("P302").Value = 7000
$ws.Range("Q302").Value = '$/bandeja 18 kilos'
$ws.Range("R302").Value = $origen
$ws.Range("S302").Value = 389
$ws.Range("T302").Value = 18

# --- Row 303 ---
$ws.Range("A303").Value = $mercadoId
$ws.Range("B303").Value = $mercado
$ws.Range("C303").Value = $region
$ws.Range("D303").Value = 44449
$ws.Range("E303").Value = $codreg
$ws.Range("F303").Value = $tipo
$ws.Range("G303").Value = $productoId
$ws.Range("H303").Value = $producto
$ws.Range("I303").Value = $categoriaId
$ws.Range("J303").Value = $categoria
$ws.Range("K303").Value = "Murcott"
$ws.Range("L303").Value = "Primera"
$ws.Range("M303").Value = 15
$ws.Range("N303").Value = 157000
$ws.Range("O303").Value = 157000
$ws.Range("P303").Value = 157000
$ws.Range("Q303").Value = '$/bins (450 kilos)'
$ws.Range("R303").Value = $origen
$ws.Range("S303").Value = 349
$ws.Range("T303").Value = 450

# --- Row 304 ---
$ws.Range("A304").Value = $mercadoId
$ws.Range("B304").Value = $mercado
$ws.Range("C304").Value = $region
$ws.Range("D304").Value = 44449
$ws.Range("E304").Value = $codreg
$ws.Range("F304").Value = $tipo
$ws.Range("G304").Value = $productoId
$ws.Range("H304").Value = $producto
$ws.Range("I304").Value = $categoriaId
$ws.Range("J304").Value = $categoria
$ws.Range("K304").Value = "Murcott"
$ws.Range("L304").Value = "Segunda"
$ws.Range("M304").Value = 6
$ws.Range("N304").Value = 135000
$ws.Range("O304").Value = 135000
$ws.Range("P304").Value = 135000
$ws.Range("Q304").Value = '$/bins (450 kilos)'
$ws.Range("R304").Value = $origen
$ws.Range("S304").Value = 300
$ws.Range("T304").Value = 450

# --- Row 305 ---
$ws.Range("A305").Value = $mercadoId
$ws.Range("B305").Value = $mercado
$ws.Range("C305").Value = $region
$ws.Range("D305").Value = 44449
$ws.Range("E305").Value = $codreg
$ws.Range("F305").Value = $tipo
$ws.Range("G305").Value = $productoId
$ws.Range("H305").Value = $producto
$ws.Range("I305").Value = $categoriaId
$ws.Range("J305").Value = $categoria
$ws.Range("K305").Value = "Murcott"
$ws.Range("L305").Value = "Tercera"
$ws.Range("M305").Value = 4
$ws.Range("N305").Value = 100000
$ws.Range("O305").Value = 100000
$ws.Range("P305").Value = 100000
$ws.Range("Q305").Value = '$/bins (450 kilos)'
$ws.Range("R305").Value = $origen
$ws.Range("S305").Value = 222
$ws.Range("T305").Value = 450
